$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loop")

$ws.Range("A3").Value = "TR ACTUALS TRY Dec"
$ws.Range("B3").Value = "TR ACTUALS TRY Apr"
$ws.Range("C3").Value = "Feb-18"

$ws.Range("C3").Select()
